$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item(1)   # Overview
$wsZhCn     = $wb.Worksheets.Item(2)   # zh-cn
$wsDeDe     = $wb.Worksheets.Item(3)   # de-de

$statusText = "Handed back: in sync with en-US"

$daeName = "dae40f3e-36d2-4979-815f-7582a5f37575.md"
$daeUrl  = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/14aeb557d9ff51af10a677d0dc6dfcc8b8aa067f/e2e/dae40f3e-36d2-4979-815f-7582a5f37575.md"

$zhCnHandbackFile = "dae40f3e-36d2-4979-815f-7582a5f37575.83d9fae8c153e8bbe686471dc0b14a9d5108f3dc.zh-cn.xlf"
$deDeHandbackFile = "dae40f3e-36d2-4979-815f-7582a5f37575.83d9fae8c153e8bbe686471dc0b14a9d5108f3dc.de-de.xlf"

$zhCnHandbackTime = "2016-09-06 17:44:40"
$deDeHandbackTime = "2016-09-06 17:44:48"

# ---------------------------------------------------------------------
# 1. Status column on the Overview sheet now reflects a completed
#    handback instead of "Ready for handoff".
# ---------------------------------------------------------------------
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText
$wsOverview.Range("E3").Value = $statusText
$wsOverview.Range("F3").Value = $statusText

# ---------------------------------------------------------------------
# 2. Same status text update on the per-locale sheets, plus fill in the
#    "Latest Target File" / "Latest Handback File" / "Latest Handback
#    DateTime" columns now that a handback report has been generated.
# ---------------------------------------------------------------------
$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("C3").Value = $statusText

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $daeUrl, "", "", $daeName)
$wsZhCn.Range("J2").Value = $zhCnHandbackFile
$wsZhCn.Range("K2").Value = $zhCnHandbackTime

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $daeUrl, "", "", $daeName)
$wsZhCn.Range("J3").Value = $zhCnHandbackFile
$wsZhCn.Range("K3").Value = $zhCnHandbackTime

$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("C3").Value = $statusText

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $daeUrl, "", "", $daeName)
$wsDeDe.Range("J2").Value = $deDeHandbackFile
$wsDeDe.Range("K2").Value = $deDeHandbackTime

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $daeUrl, "", "", $daeName)
$wsDeDe.Range("J3").Value = $deDeHandbackFile
$wsDeDe.Range("K3").Value = $deDeHandbackTime

# ---------------------------------------------------------------------
# 3. Widen the columns that now hold the longer status text / file
#    names so the handback report reads cleanly.
#    (ColumnWidth is specified in characters; Excel quantizes the
#    stored value, so these inputs are chosen to land on the desired
#    stored width.)
# ---------------------------------------------------------------------
$wsOverview.Columns.Item(5).ColumnWidth = 29.17    # E
$wsOverview.Columns.Item(6).ColumnWidth = 29.17    # F

$wsZhCn.Columns.Item(3).ColumnWidth = 29.17         # C
$wsZhCn.Columns.Item(9).ColumnWidth = 39.17         # I
$wsZhCn.Columns.Item(10).ColumnWidth = 39.17        # J

$wsDeDe.Columns.Item(3).ColumnWidth = 29.17         # C
$wsDeDe.Columns.Item(9).ColumnWidth = 39.17         # I
$wsDeDe.Columns.Item(10).ColumnWidth = 39.17        # J
